# Generate Report for Handoff
#
# Rebuilds the Overview / zh-cn / de-de sheets so a newly-handed-off file
# (7ce5caa9-aa24-4440-98ab-23be307f5d9e.md) and its two screenshot
# dependencies (6943078c-95dc-46be-8ed7-64f53bd8f7eb.png and
# 9f60a06e-865f-4617-bd2f-038c72122a42.png) show up alongside the
# pre-existing .localization-config row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the old hyperlinks so we can rebuild them in the right order/ids.
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = "6943078c-95dc-46be-8ed7-64f53bd8f7eb.png"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

$ov.Range("A3").Value = "7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Range("A4").Value = "9f60a06e-865f-4617-bd2f-038c72122a42.png"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

$ov.Range("A5").Value = ".localization-config"
$ov.Range("B5").Value = "Not to be localized"
$ov.Range("C5").Value = "Not to be localized"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/6943078c-95dc-46be-8ed7-64f53bd8f7eb.png", "", "", "6943078c-95dc-46be-8ed7-64f53bd8f7eb.png")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/7ce5caa9-aa24-4440-98ab-23be307f5d9e.md", "", "", "7ce5caa9-aa24-4440-98ab-23be307f5d9e.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/9f60a06e-865f-4617-bd2f-038c72122a42.png", "", "", "9f60a06e-865f-4617-bd2f-038c72122a42.png")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fba508942bd4f8fdefc2be0549e081ada544d9a8/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

# Row 2: screenshot dependency feeding the new handoff file.
$zh.Range("A2").Value = "6943078c-95dc-46be-8ed7-64f53bd8f7eb.png"
$zh.Range("B2").Value = "Ready for handoff"
$zh.Range("C2").Value = "95e4cec93fc78aaa83683e46f5c9c59b9406b156.png"
$zh.Range("D2").Value = "2016-03-10 23:17:26"
$zh.Range("E2").Value = ""
$zh.Range("F2").Value = ""
$zh.Range("G2").Value = "0001-01-01 00:00:00"
$zh.Range("H2").Value = "IsDependency"
$zh.Range("I2").Value = "e2e\7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"

# Row 3: the new handoff file itself.
$zh.Range("A3").Value = "7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "7ce5caa9-aa24-4440-98ab-23be307f5d9e.43e9d8f08bfa9eba5e4a75b2a7ec3dbca8b8946b.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-10 23:17:26"
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Include"

# Row 4: second screenshot dependency feeding the new handoff file.
$zh.Range("A4").Value = "9f60a06e-865f-4617-bd2f-038c72122a42.png"
$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("C4").Value = "16191e7f86cc68f9d468aa8e8acaac487ce9860f.png"
$zh.Range("D4").Value = "2016-03-10 23:17:26"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "IsDependency"
$zh.Range("I4").Value = "e2e\7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"

# Row 5: the pre-existing, not-localized config file, now pushed down.
$zh.Range("A5").Value = ".localization-config"
$zh.Range("B5").Value = "Not to be localized"
$zh.Range("D5").Value = "0001-01-01 00:00:00"
$zh.Range("G5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").Value = "Ignored"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/6943078c-95dc-46be-8ed7-64f53bd8f7eb.png", "", "", "6943078c-95dc-46be-8ed7-64f53bd8f7eb.png")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca3f45d96e848b7665faee564f6e40af67403ae6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/95e4cec93fc78aaa83683e46f5c9c59b9406b156.png", "", "", "95e4cec93fc78aaa83683e46f5c9c59b9406b156.png")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/7ce5caa9-aa24-4440-98ab-23be307f5d9e.md", "", "", "7ce5caa9-aa24-4440-98ab-23be307f5d9e.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca3f45d96e848b7665faee564f6e40af67403ae6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7ce5caa9-aa24-4440-98ab-23be307f5d9e.43e9d8f08bfa9eba5e4a75b2a7ec3dbca8b8946b.zh-cn.xlf", "", "", "7ce5caa9-aa24-4440-98ab-23be307f5d9e.43e9d8f08bfa9eba5e4a75b2a7ec3dbca8b8946b.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/9f60a06e-865f-4617-bd2f-038c72122a42.png", "", "", "9f60a06e-865f-4617-bd2f-038c72122a42.png")
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca3f45d96e848b7665faee564f6e40af67403ae6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/16191e7f86cc68f9d468aa8e8acaac487ce9860f.png", "", "", "16191e7f86cc68f9d468aa8e8acaac487ce9860f.png")
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fba508942bd4f8fdefc2be0549e081ada544d9a8/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

# Row 2: screenshot dependency feeding the new handoff file.
$de.Range("A2").Value = "6943078c-95dc-46be-8ed7-64f53bd8f7eb.png"
$de.Range("B2").Value = "Ready for handoff"
$de.Range("C2").Value = "95e4cec93fc78aaa83683e46f5c9c59b9406b156.png"
$de.Range("D2").Value = "2016-03-10 23:17:32"
$de.Range("E2").Value = ""
$de.Range("F2").Value = ""
$de.Range("G2").Value = "0001-01-01 00:00:00"
$de.Range("H2").Value = "IsDependency"
$de.Range("I2").Value = "e2e\7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"

# Row 3: the new handoff file itself.
$de.Range("A3").Value = "7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "7ce5caa9-aa24-4440-98ab-23be307f5d9e.43e9d8f08bfa9eba5e4a75b2a7ec3dbca8b8946b.de-de.xlf"
$de.Range("D3").Value = "2016-03-10 23:17:32"
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Include"

# Row 4: second screenshot dependency feeding the new handoff file.
$de.Range("A4").Value = "9f60a06e-865f-4617-bd2f-038c72122a42.png"
$de.Range("B4").Value = "Ready for handoff"
$de.Range("C4").Value = "16191e7f86cc68f9d468aa8e8acaac487ce9860f.png"
$de.Range("D4").Value = "2016-03-10 23:17:32"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "IsDependency"
$de.Range("I4").Value = "e2e\7ce5caa9-aa24-4440-98ab-23be307f5d9e.md"

# Row 5: the pre-existing, not-localized config file, now pushed down.
$de.Range("A5").Value = ".localization-config"
$de.Range("B5").Value = "Not to be localized"
$de.Range("D5").Value = "0001-01-01 00:00:00"
$de.Range("G5").Value = "0001-01-01 00:00:00"
$de.Range("H5").Value = "Ignored"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/6943078c-95dc-46be-8ed7-64f53bd8f7eb.png", "", "", "6943078c-95dc-46be-8ed7-64f53bd8f7eb.png")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68ed3005ed85cf360e2df45233f371696fbda516/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/95e4cec93fc78aaa83683e46f5c9c59b9406b156.png", "", "", "95e4cec93fc78aaa83683e46f5c9c59b9406b156.png")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/7ce5caa9-aa24-4440-98ab-23be307f5d9e.md", "", "", "7ce5caa9-aa24-4440-98ab-23be307f5d9e.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68ed3005ed85cf360e2df45233f371696fbda516/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7ce5caa9-aa24-4440-98ab-23be307f5d9e.43e9d8f08bfa9eba5e4a75b2a7ec3dbca8b8946b.de-de.xlf", "", "", "7ce5caa9-aa24-4440-98ab-23be307f5d9e.43e9d8f08bfa9eba5e4a75b2a7ec3dbca8b8946b.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/9283dbcc317be554183c66c84e3e08e9d077ff27/e2e/9f60a06e-865f-4617-bd2f-038c72122a42.png", "", "", "9f60a06e-865f-4617-bd2f-038c72122a42.png")
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68ed3005ed85cf360e2df45233f371696fbda516/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/16191e7f86cc68f9d468aa8e8acaac487ce9860f.png", "", "", "16191e7f86cc68f9d468aa8e8acaac487ce9860f.png")
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fba508942bd4f8fdefc2be0549e081ada544d9a8/.localization-config", "", "", ".localization-config")
